$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEL")

# Added velpatasvir (VEL) in vitro / subtype 3a NS5A RAS findings from Smith et
# al 2018 (pubmed 29425396) as three new rows appended below the existing data
# (through row 120). Formatting is copied from the last data row (120), but
# only into the columns the new rows actually use (A:I and N) so the
# generated cell styles line up with the rest of the "VEL" sheet.
$newRows = @(
    @{ Row = 121; Substitution = "30K+93H" },
    @{ Row = 122; Substitution = "30K+31M" },
    @{ Row = 123; Substitution = "30K+31M+93H" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A120:H120").Copy() | Out-Null
    $ws.Range("A${r}:H${r}").PasteSpecial(-4122) | Out-Null

    $ws.Range("I120").Copy() | Out-Null
    $ws.Range("I${r}").PasteSpecial(-4122) | Out-Null

    $ws.Range("N120").Copy() | Out-Null
    $ws.Range("N${r}").PasteSpecial(-4122) | Out-Null

    $ws.Range("A${r}").Value = "3a"
    $ws.Range("B${r}").Value = "NS5A"
    $ws.Range("C${r}").Value = $item.Substitution
    $ws.Range("D${r}").Value = "NA"
    $ws.Range("E${r}").Value = "NA"
    $ws.Range("F${r}").Value = "NA"
    $ws.Range("G${r}").Value = "VEL"
    $ws.Range("H${r}").Value = "in vitro"
    $ws.Range("I${r}").Value = ">10000"
    $ws.Range("N${r}").Value = 29425396
}

# Matches the author's final selection/scroll position after adding the rows.
$ws.Range("G122").Select() | Out-Null
